$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 45810 to 45811 for rows 2 through 43
$ws.Range("C2:C43").Value = 45811
